$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value2 = 44229
$ws.Range("M3").Value2 = 250
$ws.Range("N3").Value2 = 6500
$ws.Range("O3").Value2 = 7000
$ws.Range("P3").Value2 = 6750
$ws.Range("Q3").Value = '$/bandeja 5 kilos'
$ws.Range("R3").Value = 'Provincia de Curicó'
$ws.Range("S3").Value2 = 1350
$ws.Range("T3").Value2 = 5
$ws.Range("D4").Value2 = 44571
$ws.Range("K4").Value = 'Brooks'
$ws.Range("L4").Value = 'Segunda'
$ws.Range("M4").Value2 = 400
$ws.Range("N4").Value2 = 8500
$ws.Range("O4").Value2 = 9000
$ws.Range("P4").Value2 = 8750
$ws.Range("Q4").Value = '$/bandeja 10 kilos'
$ws.Range("R4").Value = 'Región de O''Higgins'
$ws.Range("S4").Value2 = 875
$ws.Range("T4").Value2 = 10
$ws.Range("D5").Value2 = 44917
$ws.Range("K5").Value = 'Bing'
$ws.Range("M5").Value2 = 400
$ws.Range("P5").Value2 = 5625
$ws.Range("Q5").Value = '$/bandeja 10 kilos'
$ws.Range("S5").Value2 = 562
$ws.Range("T5").Value2 = 10
$ws.Range("D6").Value2 = 44917
$ws.Range("L6").Value = 'Primera'
$ws.Range("M6").Value2 = 400
$ws.Range("N6").Value2 = 5000
$ws.Range("O6").Value2 = 6000
$ws.Range("P6").Value2 = 5500
$ws.Range("Q6").Value = '$/bandeja 10 kilos'
$ws.Range("S6").Value2 = 550
$ws.Range("T6").Value2 = 10
$ws.Range("D7").Value2 = 44568
$ws.Range("L7").Value = 'Segunda'
$ws.Range("M7").Value2 = 200
$ws.Range("N7").Value2 = 15000
$ws.Range("O7").Value2 = 16000
$ws.Range("P7").Value2 = 15500
$ws.Range("Q7").Value = '$/bandeja 12 kilos'
$ws.Range("R7").Value = 'Región de O''Higgins'
$ws.Range("S7").Value2 = 1292
$ws.Range("T7").Value2 = 12
$ws.Range("D8").Value2 = 44208
$ws.Range("K8").Value = 'Lapins'
$ws.Range("L8").Value = 'Segunda'
$ws.Range("M8").Value2 = 200
$ws.Range("N8").Value2 = 10500
$ws.Range("O8").Value2 = 11000
$ws.Range("P8").Value2 = 10750
$ws.Range("Q8").Value = '$/bandeja 12 kilos'
$ws.Range("R8").Value = 'Provincia de Curicó'
$ws.Range("S8").Value2 = 896
$ws.Range("T8").Value2 = 12
$ws.Range("D9").Value2 = 44901
$ws.Range("K9").Value = 'Bing'
$ws.Range("M9").Value2 = 500
$ws.Range("N9").Value2 = 12000
$ws.Range("O9").Value2 = 13000
$ws.Range("P9").Value2 = 12500
$ws.Range("Q9").Value = '$/caja 15 kilos'
$ws.Range("S9").Value2 = 833
$ws.Range("T9").Value2 = 15
$ws.Range("D10").Value2 = 44901
$ws.Range("K10").Value = 'Lapins'
$ws.Range("L10").Value = 'Primera'
$ws.Range("M10").Value2 = 500
$ws.Range("N10").Value2 = 12000
$ws.Range("O10").Value2 = 13000
$ws.Range("P10").Value2 = 12500
$ws.Range("Q10").Value = '$/caja 15 kilos'
$ws.Range("S10").Value2 = 833
$ws.Range("T10").Value2 = 15
$ws.Range("D11").Value2 = 44908
$ws.Range("K11").Value = 'Rainier'
$ws.Range("L11").Value = 'Segunda'
$ws.Range("N11").Value2 = 15000
$ws.Range("O11").Value2 = 16000
$ws.Range("P11").Value2 = 15600
$ws.Range("Q11").Value = '$/caja 10 kilos'
$ws.Range("R11").Value = 'Región de O''Higgins'
$ws.Range("S11").Value2 = 1560
$ws.Range("D12").Value2 = 44537
$ws.Range("K12").Value = 'Brooks'
$ws.Range("N12").Value2 = 29000
$ws.Range("O12").Value2 = 30000
$ws.Range("P12").Value2 = 29500
$ws.Range("Q12").Value = '$/caja 20 kilos'
$ws.Range("R12").Value = 'Región de O''Higgins'
$ws.Range("S12").Value2 = 1475
$ws.Range("T12").Value2 = 20
$ws.Range("D13").Value2 = 44161
$ws.Range("K13").Value = 'Bing'
$ws.Range("L13").Value = 'Primera'
$ws.Range("M13").Value2 = 160
$ws.Range("N13").Value2 = 39000
$ws.Range("O13").Value2 = 40000
$ws.Range("P13").Value2 = 39500
$ws.Range("Q13").Value = '$/caja 20 kilos'
$ws.Range("R13").Value = 'Provincia de Curicó'
$ws.Range("S13").Value2 = 1975
$ws.Range("T13").Value2 = 20
$ws.Range("D14").Value2 = 44914
$ws.Range("K14").Value = 'Brooks'
$ws.Range("L14").Value = 'Primera'
$ws.Range("M14").Value2 = 700
$ws.Range("N14").Value2 = 7000
$ws.Range("O14").Value2 = 8000
$ws.Range("P14").Value2 = 7429
$ws.Range("Q14").Value = '$/bandeja 10 kilos'
$ws.Range("S14").Value2 = 743
$ws.Range("T14").Value2 = 10
$ws.Range("D15").Value2 = 44914
$ws.Range("K15").Value = 'Lapins'
$ws.Range("L15").Value = 'Primera'
$ws.Range("M15").Value2 = 550
$ws.Range("P15").Value2 = 7455
$ws.Range("S15").Value2 = 746
$ws.Range("D16").Value2 = 44580
$ws.Range("K16").Value = 'Sweet Heart'
$ws.Range("M16").Value2 = 300
$ws.Range("N16").Value2 = 7000
$ws.Range("O16").Value2 = 8000
$ws.Range("P16").Value2 = 7500
$ws.Range("Q16").Value = '$/bandeja 10 kilos'
$ws.Range("R16").Value = 'Región de O''Higgins'
$ws.Range("S16").Value2 = 750
$ws.Range("T16").Value2 = 10
$ws.Range("D17").Value2 = 44594
$ws.Range("K17").Value = 'Santina'
$ws.Range("M17").Value2 = 160
$ws.Range("P17").Value2 = 5500
$ws.Range("Q17").Value = '$/bandeja 5 kilos'
$ws.Range("S17").Value2 = 1100
$ws.Range("T17").Value2 = 5
$ws.Range("D18").Value2 = 44922
$ws.Range("K18").Value = 'Bing'
$ws.Range("M18").Value2 = 200
$ws.Range("R18").Value = 'Región del Maule'
$ws.Range("D20").Value2 = 44943
$ws.Range("K20").Value = 'Santina'
$ws.Range("M20").Value2 = 600
$ws.Range("N20").Value2 = 14000
$ws.Range("O20").Value2 = 15000
$ws.Range("P20").Value2 = 14333
$ws.Range("Q20").Value = '$/caja 15 kilos'
$ws.Range("R20").Value = 'Región del Maule'
$ws.Range("S20").Value2 = 956
$ws.Range("T20").Value2 = 15
$ws.Range("D21").Value2 = 44175
$ws.Range("M21").Value2 = 270
$ws.Range("N21").Value2 = 25000
$ws.Range("O21").Value2 = 26000
$ws.Range("P21").Value2 = 25500
$ws.Range("Q21").Value = '$/caja 18 kilos'
$ws.Range("S21").Value2 = 1417
$ws.Range("T21").Value2 = 18
$ws.Range("D22").Value2 = 44921
$ws.Range("K22").Value = 'Bing'
$ws.Range("L22").Value = 'Primera'
$ws.Range("M22").Value2 = 320
$ws.Range("N22").Value2 = 7500
$ws.Range("O22").Value2 = 8000
$ws.Range("P22").Value2 = 7781
$ws.Range("S22").Value2 = 778
$ws.Range("D23").Value2 = 44931
$ws.Range("K23").Value = 'Lapins'
$ws.Range("L23").Value = 'Segunda'
$ws.Range("M23").Value2 = 250
$ws.Range("N23").Value2 = 6000
$ws.Range("O23").Value2 = 6500
$ws.Range("P23").Value2 = 6250
$ws.Range("S23").Value2 = 625
$ws.Range("D24").Value2 = 44931
$ws.Range("K24").Value = 'Lapins'
$ws.Range("L24").Value = 'Segunda'
$ws.Range("M24").Value2 = 400
$ws.Range("N24").Value2 = 3000
$ws.Range("O24").Value2 = 3300
$ws.Range("P24").Value2 = 3150
$ws.Range("Q24").Value = '$/bandeja 5 kilos'
$ws.Range("S24").Value2 = 630
$ws.Range("T24").Value2 = 5
$ws.Range("D25").Value2 = 44557
$ws.Range("M25").Value2 = 250
$ws.Range("N25").Value2 = 9000
$ws.Range("O25").Value2 = 10000
$ws.Range("P25").Value2 = 9500
$ws.Range("Q25").Value = '$/bandeja 10 kilos'
$ws.Range("R25").Value = 'Provincia de Curicó'
$ws.Range("S25").Value2 = 950
$ws.Range("T25").Value2 = 10
